# This script updates a set of odds values on the active worksheet (Sheet1)
# of the "Jogos da Semana FlashScore" workbook, matching the commit's
# XML diff (numeric <v> value changes only - no structural/type changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new numeric value
$updates = @{
    "I2"  = 3
    "K2"  = 1.83
    "M2"  = 1.14
    "N2"  = 5.5
    "Q2"  = 2.87
    "R2"  = 1.37
    "V2"  = 1.57
    "Z2"  = 29
    "AJ2" = 12
    "AN2" = 4.5
    "AX2" = 19
    "AY2" = 34

    "G3"  = 2.45
    "I3"  = 2.88
    "J3"  = 3.2
    "L3"  = 3.6
    "Q3"  = 2.2
    "R3"  = 1.62
    "U3"  = 1.83
    "V3"  = 1.83
    "W3"  = 7.5
    "Y3"  = 10
    "Z3"  = 23
    "AH3" = 8
    "AI3" = 13
    "AK3" = 29
    "AM3" = 34
    "AO3" = 15
    "AQ3" = 51
    "AW3" = 4.75
    "AX3" = 17
    "BB3" = 201

    "G4"  = 1.96

    "H5"  = 3.4
    "I5"  = 2.18
    "L5"  = 2.72
    "O5"  = 1.23
    "P5"  = 3.35
    "Q5"  = 1.7
    "R5"  = 1.93
    "U5"  = 1.57
    "V5"  = 2.1
    "AC5" = 11.75
    "AD5" = 6.7
    "AH5" = 9
    "AI5" = 11.75
    "AM5" = 23
    "AR5" = 100
    "AT5" = 2.7
    "AX5" = 11
    "AY5" = 17.5
    "AZ5" = 40
    "BA5" = 65
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
